# Applies the scheduled market-data refresh captured in the commit diff:
# updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 105.333336
$ws.Range("I8").Value = 105.333336
$ws.Range("K8").Value = 316.000008
$ws.Range("M8").Value = -177.000008
$ws.Range("H40").Value = 2263.158
$ws.Range("J40").Value = 2277.7778
$ws.Range("L40").Value = 2277.7778
$ws.Range("N40").Value = -2627.7778
$ws.Range("H42").Value = 333.33334
$ws.Range("J42").Value = 380
$ws.Range("L42").Value = 1140
$ws.Range("N42").Value = -1600
$ws.Range("H43").Value = 1559.25
$ws.Range("I43").Value = 2305.3333
$ws.Range("J43").Value = 600
$ws.Range("K43").Value = 2305.3333
$ws.Range("L43").Value = 600
$ws.Range("M43").Value = -2236.3333
$ws.Range("N43").Value = -738
$ws.Range("H52").Value = 2666.6667
$ws.Range("I52").Value = 500
$ws.Range("J52").Value = 3100
$ws.Range("K52").Value = 1500
$ws.Range("L52").Value = 9300
$ws.Range("M52").Value = -1340
$ws.Range("N52").Value = -9620
$ws.Range("H74").Value = 3640705.8
$ws.Range("I74").Value = 3920375.5
$ws.Range("K74").Value = 3920375.5
$ws.Range("M74").Value = -3919439.5
$ws.Range("H76").Value = 52383760
$ws.Range("I76").Value = 52383760
$ws.Range("K76").Value = 52383760
$ws.Range("M76").Value = -52383445
$ws.Range("H77").Value = 3640705.8
$ws.Range("I77").Value = 3920375.5
$ws.Range("K77").Value = 19601877.5
$ws.Range("M77").Value = -19597197.5
$ws.Range("H79").Value = 52383760
$ws.Range("I79").Value = 52383760
$ws.Range("K79").Value = 52383760
$ws.Range("M79").Value = -52382668
$ws.Range("H86").Value = 111114690
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 125003896
$ws.Range("K86").Value = 1000
$ws.Range("L86").Value = 125003896
$ws.Range("M86").Value = 123
$ws.Range("N86").Value = -125006142
$ws.Range("H89").Value = 111114690
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 125003896
$ws.Range("K89").Value = 5000
$ws.Range("L89").Value = 625019480
$ws.Range("M89").Value = 616
$ws.Range("N89").Value = -625030712
$ws.Range("H92").Value = 1184.2646
$ws.Range("I92").Value = 1084.2858
$ws.Range("J92").Value = 1650.8334
$ws.Range("K92").Value = 1084.2858
$ws.Range("L92").Value = 1650.8334
$ws.Range("M92").Value = 163.7141999999999
$ws.Range("N92").Value = -4146.8334
$ws.Range("H138").Value = 2966.38
$ws.Range("I138").Value = 700.13794
$ws.Range("J138").Value = 3892.028
$ws.Range("K138").Value = 2100.41382
$ws.Range("L138").Value = 11676.084
$ws.Range("M138").Value = 3039.58618
$ws.Range("N138").Value = -21956.084

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2007.762
$ws.Range("I45").Value = 1368.4117
$ws.Range("J45").Value = 4725
$ws.Range("K45").Value = 1368.4117
$ws.Range("L45").Value = 4725
$ws.Range("M45").Value = -991.4117000000001
$ws.Range("N45").Value = -5479
$ws.Range("H61").Value = 774.08
$ws.Range("I61").Value = 733.5454999999999
$ws.Range("J61").Value = 1071.3334
$ws.Range("K61").Value = 733.5454999999999
$ws.Range("L61").Value = 1071.3334
$ws.Range("M61").Value = -521.5454999999999
$ws.Range("N61").Value = -1495.3334
$ws.Range("H122").Value = 998.2857
$ws.Range("I122").Value = 998.2857
$ws.Range("K122").Value = 2994.8571
$ws.Range("M122").Value = -544.8571000000002
$ws.Range("H132").Value = 1356.8572
$ws.Range("I132").Value = 851
$ws.Range("J132").Value = 2368.5715
$ws.Range("K132").Value = 2553
$ws.Range("L132").Value = 7105.7145
$ws.Range("M132").Value = -23
$ws.Range("N132").Value = -12165.7145
$ws.Range("H136").Value = 774.08
$ws.Range("I136").Value = 733.5454999999999
$ws.Range("J136").Value = 1071.3334
$ws.Range("K136").Value = 2200.6365
$ws.Range("L136").Value = 3214.0002
$ws.Range("M136").Value = 349.3635000000004
$ws.Range("N136").Value = -8314.0002

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8440.691999999999
$ws.Range("I105").Value = 11633.167
$ws.Range("J105").Value = 5704.2856
$ws.Range("K105").Value = 11633.167
$ws.Range("L105").Value = 5704.2856
$ws.Range("M105").Value = -9886.166999999999
$ws.Range("N105").Value = -9198.285599999999
$ws.Range("H134").Value = 120631.88
$ws.Range("I134").Value = 3394
$ws.Range("J134").Value = 402002.8
$ws.Range("K134").Value = 10182
$ws.Range("L134").Value = 1206008.4
$ws.Range("M134").Value = -7647
$ws.Range("N134").Value = -1211078.4

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 28355.2
$ws.Range("J100").Value = 28355.2
$ws.Range("L100").Value = 28355.2
$ws.Range("N100").Value = -30519.2

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6985.7144
$ws.Range("I3").Value = 3200
$ws.Range("K3").Value = 9600
$ws.Range("M3").Value = -9488
$ws.Range("H113").Value = 824.96295
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 824.96295
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2474.88885
$ws.Range("N113").Value = -6814.888849999999
$ws.Range("M113").ClearContents()
$ws.Range("H131").Value = 791.2
$ws.Range("J131").Value = 810.8421
$ws.Range("L131").Value = 2432.5263
$ws.Range("N131").Value = -12512.5263

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2844.7144
$ws.Range("I102").Value = 1500
$ws.Range("J102").Value = 3382.6
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 3382.6
$ws.Range("M102").Value = 122
$ws.Range("N102").Value = -6626.6
$ws.Range("H132").Value = 2759.3513
$ws.Range("I132").Value = 2531.92
$ws.Range("J132").Value = 3233.1667
$ws.Range("K132").Value = 7595.76
$ws.Range("L132").Value = 9699.500100000001
$ws.Range("M132").Value = -5065.76
$ws.Range("N132").Value = -14759.5001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1004.2222
$ws.Range("I22").Value = 399.75
$ws.Range("J22").Value = 1258.7368
$ws.Range("K22").Value = 399.75
$ws.Range("L22").Value = 1258.7368
$ws.Range("M22").Value = -104.75
$ws.Range("N22").Value = -1848.7368
$ws.Range("H27").Value = 1004.2222
$ws.Range("I27").Value = 399.75
$ws.Range("J27").Value = 1258.7368
$ws.Range("K27").Value = 399.75
$ws.Range("L27").Value = 1258.7368
$ws.Range("M27").Value = -292.75
$ws.Range("N27").Value = -1472.7368

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1534.75
$ws.Range("I132").Value = 929.2381
$ws.Range("J132").Value = 2690.7273
$ws.Range("K132").Value = 2787.7143
$ws.Range("L132").Value = 8072.1819
$ws.Range("M132").Value = -257.7143000000001
$ws.Range("N132").Value = -13132.1819
